$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(10.56763069546768,4.236652108051314,0,18.72659692017758,42.43465696815652,37.98695560698908,16.58317817649886,0,8.66396323277201,9.9506573560273,0,17.84820150238873,20.29459104553)
    3 = @(10.34047401061864,4.050916498618962,0,18.64269901700582,42.34753485456569,38.00591867616666,16.62830501849678,0,8.681937061425256,9.806429210274601,0,17.75609966239461,20.36054486379675)
    4 = @(10.20110140270072,3.931520962450746,0,18.59472313278851,42.30435450494544,38.02857542752347,16.65875794916247,0,8.693524692733629,9.719081656557593,0,17.70283597219475,20.4029300596753)
    5 = @(10.14441931257148,3.881555914051275,0,18.57607650915467,42.28936236007718,38.0405705759162,16.6718572093637,0,8.698385876187157,9.683836004845874,0,17.68197434148266,20.4206787619059)
    6 = @(10.13501674047457,3.873181265564088,0,18.57303524322086,42.28703048232643,38.04272895920612,16.67407394909639,0,8.699201487387244,9.678005891646652,0,17.67856172272747,20.42365473192518)
    7 = @(10.20033639196061,3.930852370778092,0,18.59446797982351,42.30414175886321,38.02872602438387,16.65893181990292,0,8.693589688469745,9.718604848225281,0,17.70255118683057,20.40316749426648)
    8 = @(10.48933768862761,4.173741217521377,0,18.69694360770619,42.40248410139932,37.99120452398581,16.59816775576523,0,8.670046371168064,9.900703734398714,0,17.81577316972126,20.31694059623657)
    9 = @(11.05315429000948,4.606271930056406,0,18.92528731468651,42.67657224246092,38.00525021920998,16.50082308014222,0,8.628235655116208,10.26539833155499,0,18.06307792330862,20.16277935181317)
    10 = @(11.46084310514411,4.895953982972758,0,19.10876597895138,42.92652176250393,38.06920249075064,16.44264800289298,0,8.600146657399399,10.53531248430165,0,18.25903333562719,20.058531140353)
    11 = @(11.64395625395883,5.021416862081947,0,19.19541767453644,43.05052462309826,38.10994795007716,16.41908955238706,0,8.587933379467168,10.65797991861503,0,18.35101107785629,20.01304463359658)
    12 = @(11.71288721003236,5.068004074574109,0,19.22866787341689,43.09893702569592,38.12704997845528,16.41058714769049,0,8.583389276403008,10.70436930120746,0,18.38622510464071,19.99609718757726)
    13 = @(11.6980610364465,5.05801193577221,0,19.22148775281246,43.08844626380517,38.12329243221357,16.41239965933564,0,8.584364342725008,10.69438198377289,0,18.37862443966689,19.99973480838237)
    14 = @(11.64963586155476,5.025268190542792,0,19.19814455218641,43.0544785322828,38.11132144777891,16.41838165903031,0,8.587557916188175,10.66179789941613,0,18.35390059896283,20.01164480598619)
    15 = @(11.61991857342112,5.005091125795247,0,19.18390244952947,43.0338610206001,38.10420656982114,16.4221003538129,0,8.589524583464515,10.64182980216154,0,18.33880582860698,20.01897610265863)
    16 = @(11.44882269916593,4.887626289742912,0,19.10316541707311,42.91862298045049,38.06677395238177,16.44424614508758,0,8.600956151335604,10.52728961438551,0,18.2530774303635,20.06154265866496)
    17 = @(11.34320678040427,4.813936605262485,0,19.05443718668941,42.85054848751835,38.04679335904861,16.45857669225168,0,8.608113370095786,10.45696015860751,0,18.20119560910282,20.08815102112803)
    18 = @(11.28224247797654,4.770958548766374,0,19.02671122487831,42.81236498200735,38.03639822449526,16.46709269291943,0,8.612283173528688,10.4165000420793,0,18.17162302906847,20.10363779119762)
    19 = @(11.26156606259135,4.75630544479579,0,19.01737605724631,42.79960424747966,38.03306711303014,16.47002300351361,0,8.613704139359493,10.40280088540861,0,18.16165708768901,20.10891269885758)
    20 = @(11.35447275603795,4.82184255933365,0,19.05959336456419,42.85769480379472,38.04880678995175,16.45702287390581,0,8.607345972795221,10.46444805787214,0,18.20669090246533,20.08529965207325)
    21 = @(11.66387119210018,5.034910972657443,0,19.20498932569056,43.06441640187726,38.11479225691193,16.41661323167887,0,8.58661769637502,10.67137067331393,0,18.36115236050236,20.00813903416048)
    22 = @(11.86366074098703,5.168779361020928,0,19.30255134443902,43.207989404251,38.16766392930032,16.39264387980185,0,8.573541360414655,10.80622463981607,0,18.46432890941061,19.95932593653932)
    23 = @(11.75727350185375,5.09782813046255,0,19.25025565255299,43.13059605418729,38.1385550850589,16.40521318495849,0,8.58047749471357,10.73430002317303,0,18.409065907835,19.98523093803589)
    24 = @(11.34938016381703,4.818270185368125,0,19.05726135615077,42.85446098233917,38.04789311610735,16.45772449181074,0,8.607692741687654,10.46106285853819,0,18.20420568430599,20.08658816633501)
    25 = @(10.90144573888064,4.494109085135824,0,18.8606794575945,42.59381487306896,37.9920399051713,16.52481724718675,0,8.639082868592283,10.16620175520685,0,17.99358408716206,20.20289483011966)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    $arr = New-Object "object[,]" 1,13
    for ($i = 0; $i -lt 13; $i++) {
        $arr[0, $i] = $rowVals[$i]
    }
    $ws.Range("B$r`:N$r").Value = $arr
}
